$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data
$ws.Range("A2").Value = "B101"
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = [char]0x201C + "proyector, computadoras" + [char]0x201D
$ws.Range("F2").Value = "Anasagasti 1"

# Row 3 data
$ws.Range("A3").Value = "B102"
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = [char]0x201C + "proyector" + [char]0x201D
$ws.Range("F3").Value = "Anasagasti 1"

# Column widths (approximate, engine snaps to 1/6 character grid)
$ws.Columns.Item(3).ColumnWidth = 22.333333333333332
$ws.Columns.Item(4).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 13.0

# Create the extra (unused) underline font / style, as seen in the diff
$ws.Range("Z1").Font.Underline = 2
$ws.Range("Z1").Clear()

# Move selection to F4 (matches diff's final selection)
[void]$ws.Range("F4").Select()
